$d = $word.ActiveDocument

# 1. Remove the old "_GoBack" bookmark that currently sits after the
#    "Match names of the file with the name of the controller/service/view"
#    paragraph.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. Find the two trailing empty paragraphs (there are three empty
#    paragraphs right before the section break; the first two get text,
#    the last one stays empty) and fill them in.
$count = $d.Paragraphs.Count
$pDemo = $d.Paragraphs.Item($count - 2)
$pLink = $d.Paragraphs.Item($count - 1)

# Use a placeholder character so we can drop a *collapsed* bookmark at the
# exact end of the new sentence (immediately before the paragraph mark),
# matching the way Word leaves its auto "_GoBack" bookmark.
$pDemo.Range.Text = "A small demo with all the files can be found here:X"
$pLink.Range.Text = "https://github.com/Travo100/application-frameworks"

$pDemo = $d.Paragraphs.Item($count - 2)
$placeholderStart = $pDemo.Range.End - 2
$placeholderRange = $d.Range($placeholderStart, $placeholderStart + 1)
$d.Bookmarks.Add("_GoBack", $placeholderRange)

$deleteRange = $d.Range($placeholderStart, $placeholderStart + 1)
$deleteRange.Text = ""
